$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.94
$ws.Range("G2").Value = 1.96
$ws.Range("I2").Value = 3.85
$ws.Range("N2").Value = 6.8
$ws.Range("O2").Value = 1.15
$ws.Range("P2").Value = 2.92
$ws.Range("Q2").Value = 1.48
$ws.Range("R2").Value = 1.79
$ws.Range("S2").Value = 2.2
$ws.Range("W2").Value = 2.04
$ws.Range("X2").Value = 29
$ws.Range("Y2").Value = 24
$ws.Range("AB2").Value = 16.5
$ws.Range("AF2").Value = 16.5
$ws.Range("AN2").Value = 7.6
# Row 3
$ws.Range("N3").Value = 5.1
$ws.Range("O3").Value = 1.08
# Row 4
$ws.Range("F4").Value = 4.7
$ws.Range("G4").Value = 5.3
$ws.Range("I4").Value = 1.76
# Row 5
$ws.Range("H5").Value = 4.1
$ws.Range("J5").Value = 2.68
# Row 6
$ws.Range("F6").Value = 1.39
$ws.Range("G6").Value = 1.46
$ws.Range("H6").Value = 9.199999999999999
$ws.Range("J6").Value = 4.1
$ws.Range("N6").Value = 3.05
$ws.Range("O6").Value = 1.39
$ws.Range("P6").Value = 1.7
$ws.Range("W6").Value = 3.15
$ws.Range("AB6").Value = 1000
$ws.Range("AN6").Value = 1000
# Row 8
$ws.Range("I8").Value = 11.5
$ws.Range("N8").Value = 5.9
$ws.Range("P8").Value = 2.66
$ws.Range("R8").Value = 1.69
$ws.Range("U8").Value = 1.96
$ws.Range("W8").Value = 3.5
# Row 9
$ws.Range("H9").Value = 2.18
$ws.Range("K9").Value = 3.95
$ws.Range("N9").Value = 3.75
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 1.91
$ws.Range("Q9").Value = 1.9
$ws.Range("R9").Value = 1.36
$ws.Range("T9").Value = 1.72
$ws.Range("U9").Value = 2.12
$ws.Range("AJ9").Value = 65
# Row 10
$ws.Range("F10").Value = 2.36
$ws.Range("H10").Value = 2.98
$ws.Range("I10").Value = 3.15
$ws.Range("J10").Value = 3.65
$ws.Range("K10").Value = 3.9
$ws.Range("N10").Value = 4.1
$ws.Range("P10").Value = 2.1
$ws.Range("Q10").Value = 1.75
$ws.Range("R10").Value = 1.37
$ws.Range("S10").Value = 2.66
$ws.Range("X10").Value = 18.5
$ws.Range("Y10").Value = 17.5
$ws.Range("Z10").Value = 28
$ws.Range("AA10").Value = 65
$ws.Range("AB10").Value = 14.5
$ws.Range("AC10").Value = 10.5
$ws.Range("AD10").Value = 16.5
$ws.Range("AE10").Value = 42
$ws.Range("AF10").Value = 20
$ws.Range("AG10").Value = 14
$ws.Range("AH10").Value = 20
$ws.Range("AI10").Value = 48
$ws.Range("AJ10").Value = 38
$ws.Range("AK10").Value = 29
$ws.Range("AL10").Value = 42
$ws.Range("AM10").Value = 90
$ws.Range("AN10").Value = 19.5
$ws.Range("AO10").Value = 34
# Row 11
$ws.Range("F11").Value = 2.08
$ws.Range("G11").Value = 2.2
$ws.Range("H11").Value = 3.15
$ws.Range("I11").Value = 3.6
$ws.Range("K11").Value = 4.6
$ws.Range("L11").Value = 1.18
$ws.Range("N11").Value = 5.9
$ws.Range("O11").Value = 1.16
$ws.Range("P11").Value = 2.7
$ws.Range("Q11").Value = 1.48
$ws.Range("R11").Value = 1.71
$ws.Range("S11").Value = 2.16
$ws.Range("T11").Value = 1.48
$ws.Range("U11").Value = 2.68
$ws.Range("V11").Value = 1.39
$ws.Range("W11").Value = 1.83
$ws.Range("X11").Value = 32
$ws.Range("Z11").Value = 980
$ws.Range("AA11").Value = 60
$ws.Range("AF11").Value = 19
$ws.Range("AI11").Value = 36
$ws.Range("AM11").Value = 55
$ws.Range("AN11").Value = 9.4
$ws.Range("AO11").Value = 23
# Row 12
$ws.Range("S12").Value = 3.65
$ws.Range("AD12").Value = 1000
# Row 13
$ws.Range("F13").Value = 7
$ws.Range("I13").Value = 1.49
$ws.Range("J13").Value = 4.9
$ws.Range("K13").Value = 5.8
$ws.Range("N13").Value = 5.4
$ws.Range("P13").Value = 2.5
$ws.Range("Q13").Value = 1.55
$ws.Range("R13").Value = 1.61
$ws.Range("S13").Value = 2.34
$ws.Range("U13").Value = 2.1
$ws.Range("AI13").Value = 1000
# Row 14
$ws.Range("F14").Value = 2.18
$ws.Range("H14").Value = 3.25
$ws.Range("J14").Value = 3.6
$ws.Range("K14").Value = 4.2
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 4.1
$ws.Range("O14").Value = 1.29
$ws.Range("P14").Value = 2.08
$ws.Range("Q14").Value = 1.78
$ws.Range("R14").Value = 1.42
$ws.Range("S14").Value = 2.96
$ws.Range("T14").Value = 1.64
$ws.Range("U14").Value = 2.24
$ws.Range("W14").Value = 1.73
$ws.Range("X14").Value = 22
$ws.Range("AM14").Value = 100
# Row 15
$ws.Range("H15").Value = 11.5
$ws.Range("Q15").Value = 1.32
$ws.Range("U15").Value = 1.98
# Row 16
$ws.Range("G16").Value = 2.28
$ws.Range("H16").Value = 3.3
$ws.Range("I16").Value = 3.95
$ws.Range("J16").Value = 3.65
$ws.Range("Q16").Value = 1.76
$ws.Range("R16").Value = 1.42
$ws.Range("S16").Value = 2.94
$ws.Range("T16").Value = 1.66
$ws.Range("W16").Value = 1.78
$ws.Range("AD16").Value = 990
$ws.Range("AF16").Value = 18
$ws.Range("AI16").Value = 55
$ws.Range("AJ16").Value = 34
# Row 17
$ws.Range("F17").Value = 5.5
$ws.Range("H17").Value = 1.6
$ws.Range("I17").Value = 1.71
$ws.Range("J17").Value = 4
$ws.Range("K17").Value = 4.7
# Row 20
$ws.Range("I20").Value = 1.57
$ws.Range("T20").Value = 1.62
$ws.Range("V20").Value = 2.74
# Row 21
$ws.Range("F21").Value = 5.2
$ws.Range("I21").Value = 1.8
$ws.Range("L21").Value = 1.32
$ws.Range("N21").Value = 3.8
$ws.Range("P21").Value = 1.98
$ws.Range("Q21").Value = 1.86
$ws.Range("R21").Value = 1.38
$ws.Range("S21").Value = 3.1
$ws.Range("V21").Value = 2.24
$ws.Range("AE21").Value = 22
$ws.Range("AH21").Value = 24
$ws.Range("AI21").Value = 42
$ws.Range("AK21").Value = 90
$ws.Range("AL21").Value = 75
$ws.Range("AO21").Value = 13
# Row 22
$ws.Range("F22").Value = 2.42
$ws.Range("G22").Value = 2.56
$ws.Range("H22").Value = 3.45
$ws.Range("I22").Value = 3.7
$ws.Range("P22").Value = 1.62
$ws.Range("Q22").Value = 2.48
# Row 23
$ws.Range("M23").Value = 1.03
$ws.Range("P23").Value = 2.42
$ws.Range("S23").Value = 2.42
$ws.Range("U23").Value = 2.48
$ws.Range("AJ23").Value = 34
$ws.Range("AK23").Value = 25
$ws.Range("AO23").Value = 24
# Row 24
$ws.Range("G24").Value = 1.24
$ws.Range("H24").Value = 13.5
$ws.Range("J24").Value = 8.199999999999999
$ws.Range("K24").Value = 9.4
$ws.Range("N24").Value = 1.1
$ws.Range("Q24").Value = 1.27
$ws.Range("T24").Value = 1.04
$ws.Range("U24").Value = 1.04
$ws.Range("AF24").Value = 12.5
$ws.Range("AJ24").Value = 14
$ws.Range("AK24").Value = 13.5
$ws.Range("AN24").Value = 2.84
# Row 25
$ws.Range("H25").Value = 2.78
$ws.Range("L25").Value = 1.32
$ws.Range("N25").Value = 5.1
$ws.Range("O25").Value = 1.23
$ws.Range("P25").Value = 2.4
$ws.Range("Q25").Value = 1.7
$ws.Range("S25").Value = 2.68
$ws.Range("U25").Value = 2.72
$ws.Range("Y25").Value = 14.5
$ws.Range("AB25").Value = 14.5
$ws.Range("AF25").Value = 21
$ws.Range("AN25").Value = 16.5
$ws.Range("AO25").Value = 18
